$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.536300000000003
$ws.Range("B14").Value = 5.4232
$ws.Range("D15").Value = -8.005999999999998
$ws.Range("B16").Value = 6.033099999999999
$ws.Range("B21").Value = 8.220900000000002
$ws.Range("D21").Value = -8.34
$ws.Range("D22").Value = -8.058600000000002
$ws.Range("B23").Value = 8.981899999999992
$ws.Range("D24").Value = -7.394099999999997
$ws.Range("B25").Value = 5.376100000000003
$ws.Range("B26").Value = 5.459100000000007
$ws.Range("D27").Value = -8.503199999999998
$ws.Range("D28").Value = -8.330799999999996
$ws.Range("B29").Value = 5.046500000000005
$ws.Range("D36").Value = -6.540899999999999
$ws.Range("D39").Value = -7.917699999999996
$ws.Range("B40").Value = 8.982899999999999
$ws.Range("D45").Value = -7.222800000000001
$ws.Range("D48").Value = -7.455199999999996
$ws.Range("D49").Value = -7.9813
$ws.Range("D52").Value = -8.079800000000006
$ws.Range("B53").Value = 5.784400000000004
$ws.Range("D53").Value = -8.569400000000003
$ws.Range("D54").Value = -8.675999999999998
$ws.Range("B57").Value = 4.603199999999996
$ws.Range("D57").Value = -7.939499999999998
$ws.Range("B59").Value = 5.573300000000004
$ws.Range("B65").Value = 5.5344
$ws.Range("B69").Value = 5.381699999999995
$ws.Range("D70").Value = -7.095499999999999
$ws.Range("D71").Value = -7.259099999999997
$ws.Range("B79").Value = 9.772300000000001
$ws.Range("B83").Value = 4.963499999999996
$ws.Range("D86").Value = -7.848699999999993
$ws.Range("D87").Value = -8.207599999999998
$ws.Range("D89").Value = -8.263100000000003
$ws.Range("B91").Value = 5.540600000000004
$ws.Range("B93").Value = 5.329600000000002
$ws.Range("B100").Value = 5.619700000000003
$ws.Range("D101").Value = -8.244999999999999
$ws.Range("B103").Value = 5.764399999999995
